$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Title and header text: October -> November ---
$ws.Range("A1").Value = "Table 4.6.B. Receipts of Coal Delivered for Electricity Generation by State, (Year-to-Date) November 2016 and 2015"
$ws.Range("B5").Value = "November 2016 YTD"
$ws.Range("E5").Value = "November 2016 YTD"
$ws.Range("G5").Value = "November 2016 YTD"
$ws.Range("I5").Value = "November 2016 YTD"
$ws.Range("K5").Value = "November 2016 YTD"
$ws.Range("C5").Value = "November 2015 YTD"
$ws.Range("F5").Value = "November 2015 YTD"
$ws.Range("H5").Value = "November 2015 YTD"
$ws.Range("J5").Value = "November 2015 YTD"
$ws.Range("L5").Value = "November 2015 YTD"

# --- Data cell updates (rows 6-67) ---
$ws.Range("B6").Value = 1162
$ws.Range("C6").Value = 1929
$ws.Range("D6").Value = -0.4
$ws.Range("F6").Value = 657
$ws.Range("G6").Value = 971
$ws.Range("H6").Value = 1244
$ws.Range("K6").Value = 15
$ws.Range("L6").Value = 28
$ws.Range("B8").Value = 76
$ws.Range("C8").Value = 97
$ws.Range("D8").Value = -0.21
$ws.Range("G8").Value = 62
$ws.Range("H8").Value = 69
$ws.Range("K8").Value = 15
$ws.Range("L8").Value = 28
$ws.Range("B9").Value = 824
$ws.Range("C9").Value = 925
$ws.Range("D9").Value = -0.11
$ws.Range("G9").Value = 824
$ws.Range("H9").Value = 925
$ws.Range("C10").Value = 657
$ws.Range("D10").Value = -0.73
$ws.Range("F10").Value = 657
$ws.Range("B13").Value = 17999
$ws.Range("C13").Value = 29606
$ws.Range("D13").Value = -0.39
$ws.Range("G13").Value = 17648
$ws.Range("H13").Value = 28804
$ws.Range("K13").Value = 351
$ws.Range("L13").Value = 802
$ws.Range("B14").Value = 581
$ws.Range("C14").Value = 785
$ws.Range("D14").Value = -0.26
$ws.Range("G14").Value = 581
$ws.Range("H14").Value = 785
$ws.Range("B15").Value = 534
$ws.Range("C15").Value = 953
$ws.Range("D15").Value = -0.44
$ws.Range("G15").Value = 273
$ws.Range("H15").Value = 661
$ws.Range("K15").Value = 261
$ws.Range("L15").Value = 292
$ws.Range("B16").Value = 16884
$ws.Range("C16").Value = 27868
$ws.Range("D16").Value = -0.39
$ws.Range("G16").Value = 16793
$ws.Range("H16").Value = 27359
$ws.Range("L16").Value = 510
$ws.Range("B17").Value = 123959
$ws.Range("C17").Value = 160832
$ws.Range("E17").Value = 73854
$ws.Range("F17").Value = 92812
$ws.Range("G17").Value = 47979
$ws.Range("H17").Value = 65267
$ws.Range("K17").Value = 2126
$ws.Range("L17").Value = 2726
$ws.Range("B18").Value = 35082
$ws.Range("C18").Value = 50924
$ws.Range("E18").Value = 6306
$ws.Range("F18").Value = 7877
$ws.Range("G18").Value = 27201
$ws.Range("H18").Value = 41020
$ws.Range("K18").Value = 1576
$ws.Range("L18").Value = 2027
$ws.Range("B19").Value = 26341
$ws.Range("C19").Value = 32439
$ws.Range("E19").Value = 24716
$ws.Range("F19").Value = 29817
$ws.Range("G19").Value = 1624
$ws.Range("H19").Value = 2622
$ws.Range("B20").Value = 19092
$ws.Range("C20").Value = 27272
$ws.Range("D20").Value = -0.3
$ws.Range("E20").Value = 18852
$ws.Range("F20").Value = 26980
$ws.Range("G20").Value = 229
$ws.Range("H20").Value = 242
$ws.Range("L20").Value = 23
$ws.Range("B21").Value = 26271
$ws.Range("C21").Value = 28916
$ws.Range("D21").Value = -0.091
$ws.Range("E21").Value = 7195
$ws.Range("F21").Value = 7311
$ws.Range("G21").Value = 18925
$ws.Range("H21").Value = 21383
$ws.Range("L21").Value = 222
$ws.Range("B22").Value = 17173
$ws.Range("C22").Value = 21282
$ws.Range("D22").Value = -0.19
$ws.Range("E22").Value = 16785
$ws.Range("F22").Value = 20828
$ws.Range("K22").Value = 388
$ws.Range("L22").Value = 454
$ws.Range("B23").Value = 105510
$ws.Range("C23").Value = 127958
$ws.Range("E23").Value = 104353
$ws.Range("F23").Value = 124799
$ws.Range("I23").Value = 47
$ws.Range("J23").Value = 73
$ws.Range("K23").Value = 1110
$ws.Range("L23").Value = 3086
$ws.Range("B24").Value = 15991
$ws.Range("C24").Value = 20883
$ws.Range("D24").Value = -0.23
$ws.Range("E24").Value = 14881
$ws.Range("F24").Value = 18854
$ws.Range("K24").Value = 1110
$ws.Range("L24").Value = 2029
$ws.Range("B25").Value = 13201
$ws.Range("C25").Value = 16431
$ws.Range("D25").Value = -0.2
$ws.Range("E25").Value = 13201
$ws.Range("F25").Value = 16431
$ws.Range("B26").Value = 11343
$ws.Range("C26").Value = 16138
$ws.Range("D26").Value = -0.3
$ws.Range("E26").Value = 11343
$ws.Range("F26").Value = 15787
$ws.Range("L26").Value = 336
$ws.Range("B27").Value = 32091
$ws.Range("C27").Value = 37951
$ws.Range("D27").Value = -0.15
$ws.Range("E27").Value = 32044
$ws.Range("F27").Value = 37892
$ws.Range("I27").Value = 47
$ws.Range("J27").Value = 59
$ws.Range("B28").Value = 11578
$ws.Range("C28").Value = 13920
$ws.Range("E28").Value = 11578
$ws.Range("F28").Value = 13199
$ws.Range("L28").Value = 721
$ws.Range("B29").Value = 20068
$ws.Range("C29").Value = 21673
$ws.Range("D29").Value = -0.074
$ws.Range("E29").Value = 20068
$ws.Range("F29").Value = 21673
$ws.Range("B30").Value = 1238
$ws.Range("C30").Value = 962
$ws.Range("D30").Value = 0.29
$ws.Range("E30").Value = 1238
$ws.Range("F30").Value = 962
$ws.Range("B31").Value = 86575
$ws.Range("C31").Value = 103420
$ws.Range("D31").Value = -0.16
$ws.Range("E31").Value = 75346
$ws.Range("F31").Value = 89472
$ws.Range("G31").Value = 10463
$ws.Range("H31").Value = 12133
$ws.Range("K31").Value = 766
$ws.Range("L31").Value = 1816
$ws.Range("B32").Value = 214
$ws.Range("D32").Value = 0.41
$ws.Range("G32").Value = 214
$ws.Range("B34").Value = 13975
$ws.Range("C34").Value = 18049
$ws.Range("E34").Value = 13685
$ws.Range("F34").Value = 17268
$ws.Range("H34").Value = 564
$ws.Range("L34").Value = 217
$ws.Range("B35").Value = 15878
$ws.Range("C35").Value = 18727
$ws.Range("D35").Value = -0.15
$ws.Range("E35").Value = 15752
$ws.Range("F35").Value = 18542
$ws.Range("K35").Value = 126
$ws.Range("L35").Value = 185
$ws.Range("B36").Value = 4899
$ws.Range("C36").Value = 6031
$ws.Range("D36").Value = -0.19
$ws.Range("G36").Value = 4692
$ws.Range("H36").Value = 5779
$ws.Range("K36").Value = 207
$ws.Range("L36").Value = 251
$ws.Range("B37").Value = 10764
$ws.Range("C37").Value = 14940
$ws.Range("D37").Value = -0.28
$ws.Range("E37").Value = 10764
$ws.Range("F37").Value = 14375
$ws.Range("H37").Value = 220
$ws.Range("L37").Value = 345
$ws.Range("B38").Value = 7358
$ws.Range("C38").Value = 10318
$ws.Range("D38").Value = -0.29
$ws.Range("E38").Value = 7273
$ws.Range("F38").Value = 10173
$ws.Range("K38").Value = 85
$ws.Range("L38").Value = 146
$ws.Range("B39").Value = 7111
$ws.Range("C39").Value = 7170
$ws.Range("D39").Value = -0.008
$ws.Range("E39").Value = 6390
$ws.Range("F39").Value = 6050
$ws.Range("G39").Value = 443
$ws.Range("H39").Value = 805
$ws.Range("K39").Value = 279
$ws.Range("L39").Value = 315
$ws.Range("B40").Value = 26376
$ws.Range("C40").Value = 28032
$ws.Range("D40").Value = -0.059
$ws.Range("E40").Value = 21482
$ws.Range("F40").Value = 23064
$ws.Range("G40").Value = 4825
$ws.Range("H40").Value = 4612
$ws.Range("L40").Value = 357
$ws.Range("B41").Value = 61847
$ws.Range("C41").Value = 71775
$ws.Range("D41").Value = -0.14
$ws.Range("E41").Value = 57946
$ws.Range("F41").Value = 67456
$ws.Range("G41").Value = 2807
$ws.Range("H41").Value = 2995
$ws.Range("K41").Value = 1094
$ws.Range("L41").Value = 1324
$ws.Range("B42").Value = 14973
$ws.Range("C42").Value = 18845
$ws.Range("E42").Value = 14973
$ws.Range("F42").Value = 18845
$ws.Range("B43").Value = 34140
$ws.Range("C43").Value = 38021
$ws.Range("D43").Value = -0.1
$ws.Range("E43").Value = 34140
$ws.Range("F43").Value = 38021
$ws.Range("B44").Value = 4027
$ws.Range("C44").Value = 4890
$ws.Range("D44").Value = -0.18
$ws.Range("E44").Value = 1220
$ws.Range("F44").Value = 1894
$ws.Range("G44").Value = 2807
$ws.Range("H44").Value = 2995
$ws.Range("B45").Value = 8707
$ws.Range("C45").Value = 10019
$ws.Range("D45").Value = -0.13
$ws.Range("E45").Value = 7613
$ws.Range("F45").Value = 8695
$ws.Range("K45").Value = 1094
$ws.Range("L45").Value = 1324
$ws.Range("B46").Value = 98440
$ws.Range("C46").Value = 125907
$ws.Range("D46").Value = -0.22
$ws.Range("E46").Value = 47669
$ws.Range("F46").Value = 65517
$ws.Range("G46").Value = 50714
$ws.Range("H46").Value = 59828
$ws.Range("K46").Value = 57
$ws.Range("L46").Value = 562
$ws.Range("B47").Value = 11378
$ws.Range("C47").Value = 13917
$ws.Range("D47").Value = -0.18
$ws.Range("E47").Value = 9525
$ws.Range("F47").Value = 11613
$ws.Range("G47").Value = 1797
$ws.Range("H47").Value = 2228
$ws.Range("K47").Value = 57
$ws.Range("L47").Value = 75
$ws.Range("B48").Value = 6179
$ws.Range("C48").Value = 9969
$ws.Range("D48").Value = -0.38
$ws.Range("E48").Value = 4626
$ws.Range("F48").Value = 5832
$ws.Range("G48").Value = 1554
$ws.Range("H48").Value = 4137
$ws.Range("B49").Value = 9629
$ws.Range("C49").Value = 17247
$ws.Range("D49").Value = -0.44
$ws.Range("E49").Value = 8640
$ws.Range("F49").Value = 15517
$ws.Range("G49").Value = 989
$ws.Range("H49").Value = 1244
$ws.Range("L49").Value = 487
$ws.Range("B50").Value = 71254
$ws.Range("C50").Value = 84774
$ws.Range("D50").Value = -0.16
$ws.Range("E50").Value = 24879
$ws.Range("F50").Value = 32555
$ws.Range("G50").Value = 46374
$ws.Range("H50").Value = 52219
$ws.Range("B51").Value = 80806
$ws.Range("C51").Value = 97229
$ws.Range("D51").Value = -0.17
$ws.Range("E51").Value = 72339
$ws.Range("F51").Value = 86762
$ws.Range("G51").Value = 8236
$ws.Range("H51").Value = 10241
$ws.Range("B52").Value = 14066
$ws.Range("C52").Value = 20185
$ws.Range("E52").Value = 14066
$ws.Range("F52").Value = 20185
$ws.Range("B53").Value = 14591
$ws.Range("C53").Value = 16927
$ws.Range("D53").Value = -0.14
$ws.Range("E53").Value = 14591
$ws.Range("F53").Value = 16927
$ws.Range("B55").Value = 7783
$ws.Range("C55").Value = 9079
$ws.Range("D55").Value = -0.14
$ws.Range("F55").Value = 200
$ws.Range("G55").Value = 7783
$ws.Range("H55").Value = 8879
$ws.Range("B56").Value = 912
$ws.Range("C56").Value = 1238
$ws.Range("D56").Value = -0.26
$ws.Range("E56").Value = 459
$ws.Range("F56").Value = 772
$ws.Range("G56").Value = 453
$ws.Range("H56").Value = 467
$ws.Range("B57").Value = 9801
$ws.Range("C57").Value = 11285
$ws.Range("D57").Value = -0.13
$ws.Range("E57").Value = 9801
$ws.Range("F57").Value = 11285
$ws.Range("B58").Value = 12139
$ws.Range("C58").Value = 14057
$ws.Range("D58").Value = -0.14
$ws.Range("E58").Value = 11908
$ws.Range("F58").Value = 13424
$ws.Range("H58").Value = 406
$ws.Range("B59").Value = 21513
$ws.Range("C59").Value = 24458
$ws.Range("D59").Value = -0.12
$ws.Range("E59").Value = 21513
$ws.Range("F59").Value = 23969
$ws.Range("H59").Value = 489
$ws.Range("B60").Value = 4327
$ws.Range("C60").Value = 4976
$ws.Range("D60").Value = -0.13
$ws.Range("E60").Value = 951
$ws.Range("F60").Value = 1324
$ws.Range("G60").Value = 2795
$ws.Range("H60").Value = 3035
$ws.Range("K60").Value = 580
$ws.Range("L60").Value = 617
$ws.Range("B61").Value = 580
$ws.Range("C61").Value = 617
$ws.Range("D61").Value = -0.06
$ws.Range("K61").Value = 580
$ws.Range("L61").Value = 617
$ws.Range("B62").Value = 951
$ws.Range("C62").Value = 1324
$ws.Range("D62").Value = -0.28
$ws.Range("E62").Value = 951
$ws.Range("F62").Value = 1324
$ws.Range("B63").Value = 2795
$ws.Range("C63").Value = 3035
$ws.Range("D63").Value = -0.079
$ws.Range("G63").Value = 2795
$ws.Range("H63").Value = 3035
$ws.Range("B64").Value = 965
$ws.Range("C64").Value = 727
$ws.Range("D64").Value = 0.33
$ws.Range("E64").Value = 181
$ws.Range("F64").Value = 128
$ws.Range("G64").Value = 784
$ws.Range("H64").Value = 599
$ws.Range("B65").Value = 181
$ws.Range("C65").Value = 128
$ws.Range("D65").Value = 0.41
$ws.Range("E65").Value = 181
$ws.Range("F65").Value = 128
$ws.Range("B66").Value = 784
$ws.Range("C66").Value = 599
$ws.Range("D66").Value = 0.31
$ws.Range("G66").Value = 784
$ws.Range("H66").Value = 599
$ws.Range("B67").Value = 581589
$ws.Range("C67").Value = 724360
$ws.Range("D67").Value = -0.2
$ws.Range("E67").Value = 432816
$ws.Range("F67").Value = 528926
$ws.Range("G67").Value = 142397
$ws.Range("H67").Value = 184146
$ws.Range("I67").Value = 47
$ws.Range("J67").Value = 100
$ws.Range("K67").Value = 6329
$ws.Range("L67").Value = 11188
